# Update "想去人数" (want-to-go count) values in column F
# for sheets "展览" and "全部类型", per the latest generated data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3984
$ws1.Range("F4").Value = 2337
$ws1.Range("F5").Value = 465
$ws1.Range("F7").Value = 32
$ws1.Range("F11").Value = 57
$ws1.Range("F12").Value = 124
$ws1.Range("F13").Value = 1483
$ws1.Range("F15").Value = 2760

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3984
$ws4.Range("F4").Value = 2337
$ws4.Range("F5").Value = 465
$ws4.Range("F7").Value = 32
$ws4.Range("F12").Value = 57
$ws4.Range("F13").Value = 124
$ws4.Range("F16").Value = 1483
$ws4.Range("F18").Value = 2760
